$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Real Madrid"
$ws.Cells.Item(2, 2).Value = 5.28174123337364
$ws.Cells.Item(2, 3).Value = 7.183553597650514
$ws.Cells.Item(2, 4).Value = 0.6770833333333334
$ws.Cells.Item(2, 5).Value = 13
$ws.Cells.Item(2, 6).Value = 0.05932081960166213
$ws.Cells.Item(2, 7).Value = 40

$ws.Cells.Item(3, 1).Value = "Barcelona"
$ws.Cells.Item(3, 2).Value = 3.889784946236559
$ws.Cells.Item(3, 3).Value = 8.716510903426791
$ws.Cells.Item(3, 4).Value = 0.7921889568009945
$ws.Cells.Item(3, 5).Value = 12
$ws.Cells.Item(3, 6).Value = 0.04400647948164147
$ws.Cells.Item(3, 7).Value = 55

$ws.Cells.Item(4, 1).Value = "Villarreal"
$ws.Cells.Item(4, 2).Value = 5.454545454545454
$ws.Cells.Item(4, 3).Value = 7.023569023569023
$ws.Cells.Item(4, 4).Value = 0.4080135440180587
$ws.Cells.Item(4, 5).Value = 29
$ws.Cells.Item(4, 6).Value = 0.08554383722168957
$ws.Cells.Item(4, 7).Value = 23

$ws.Cells.Item(5, 1).Value = "Atlético Madrid"
$ws.Cells.Item(5, 2).Value = 4.965801886792453
$ws.Cells.Item(5, 3).Value = 9.034542314335061
$ws.Cells.Item(5, 4).Value = 0.5483838485550298
$ws.Cells.Item(5, 5).Value = 16
$ws.Cells.Item(5, 6).Value = 0.06915150539365642
$ws.Cells.Item(5, 7).Value = 24

$ws.Cells.Item(6, 1).Value = "Real Betis"
$ws.Cells.Item(6, 2).Value = 5.744583808437857
$ws.Cells.Item(6, 3).Value = 8.195079086115992
$ws.Cells.Item(6, 4).Value = 0.5107496936615796
$ws.Cells.Item(6, 5).Value = 13
$ws.Cells.Item(6, 6).Value = 0.07963232617400801
$ws.Cells.Item(6, 7).Value = 11

$ws.Cells.Item(7, 1).Value = "Espanyol"
$ws.Cells.Item(7, 2).Value = 5.886586695747001
$ws.Cells.Item(7, 3).Value = 6.951219512195122
$ws.Cells.Item(7, 4).Value = 0.3997421018697614
$ws.Cells.Item(7, 5).Value = 10
$ws.Cells.Item(7, 6).Value = 0.0955770801727156
$ws.Cells.Item(7, 7).Value = 4

$ws.Cells.Item(8, 1).Value = "Celta Vigo"
$ws.Cells.Item(8, 2).Value = 5.757462686567164
$ws.Cells.Item(8, 3).Value = 9.135042735042735
$ws.Cells.Item(8, 4).Value = 0.4984426805988144
$ws.Cells.Item(8, 5).Value = 19
$ws.Cells.Item(8, 6).Value = 0.07327046649420012
$ws.Cells.Item(8, 7).Value = 2

$ws.Cells.Item(9, 1).Value = "Real Sociedad"
$ws.Cells.Item(9, 2).Value = 4.451158106747231
$ws.Cells.Item(9, 3).Value = 7.630742049469965
$ws.Cells.Item(9, 4).Value = 0.4768534632788026
$ws.Cells.Item(9, 5).Value = 9
$ws.Cells.Item(9, 6).Value = 0.0721760391198044
$ws.Cells.Item(9, 7).Value = 9

$ws.Cells.Item(10, 1).Value = "Athletic Club"
$ws.Cells.Item(10, 2).Value = 5.917127071823204
$ws.Cells.Item(10, 3).Value = 5.264367816091954
$ws.Cells.Item(10, 4).Value = 0.5959473966720343
$ws.Cells.Item(10, 5).Value = 21
$ws.Cells.Item(10, 6).Value = 0.078125
$ws.Cells.Item(10, 7).Value = -1

$ws.Cells.Item(11, 1).Value = "Osasuna"
$ws.Cells.Item(11, 2).Value = 5.214606741573034
$ws.Cells.Item(11, 3).Value = 6.492007104795738
$ws.Cells.Item(11, 4).Value = 0.4165858389912706
$ws.Cells.Item(11, 5).Value = 11
$ws.Cells.Item(11, 6).Value = 0.09958737593398015
$ws.Cells.Item(11, 7).Value = -14

$ws.Cells.Item(12, 1).Value = "Getafe"
$ws.Cells.Item(12, 2).Value = 4.447394296951819
$ws.Cells.Item(12, 3).Value = 4.351247600767755
$ws.Cells.Item(12, 4).Value = 0.4164494642339994
$ws.Cells.Item(12, 5).Value = 14
$ws.Cells.Item(12, 6).Value = 0.1321531494442157
$ws.Cells.Item(12, 7).Value = -17

$ws.Cells.Item(13, 1).Value = "Sevilla"
$ws.Cells.Item(13, 2).Value = 3.525100401606426
$ws.Cells.Item(13, 3).Value = 7.857632933104632
$ws.Cells.Item(13, 4).Value = 0.523877001558736
$ws.Cells.Item(13, 5).Value = 19
$ws.Cells.Item(13, 6).Value = 0.08727092620108964
$ws.Cells.Item(13, 7).Value = -36

$ws.Cells.Item(14, 1).Value = "Deportivo Alavés"
$ws.Cells.Item(14, 2).Value = 4.157099697885196
$ws.Cells.Item(14, 3).Value = 5.817538896746817
$ws.Cells.Item(14, 4).Value = 0.5128296507483963
$ws.Cells.Item(14, 5).Value = 12
$ws.Cells.Item(14, 6).Value = 0.07747695852534563
$ws.Cells.Item(14, 7).Value = 1

$ws.Cells.Item(15, 1).Value = "Valencia"
$ws.Cells.Item(15, 2).Value = 5.069711538461538
$ws.Cells.Item(15, 3).Value = 7.987544483985765
$ws.Cells.Item(15, 4).Value = 0.4341463414634146
$ws.Cells.Item(15, 5).Value = 12
$ws.Cells.Item(15, 6).Value = 0.0885160253189993
$ws.Cells.Item(15, 7).Value = -12

$ws.Cells.Item(16, 1).Value = "Girona FC"
$ws.Cells.Item(16, 2).Value = 5.493638676844784
$ws.Cells.Item(16, 3).Value = 8.45362563237774
$ws.Cells.Item(16, 4).Value = 0.4401503582755785
$ws.Cells.Item(16, 5).Value = 20
$ws.Cells.Item(16, 6).Value = 0.07318728189220629
$ws.Cells.Item(16, 7).Value = -14

$ws.Cells.Item(17, 1).Value = "Elche"
$ws.Cells.Item(17, 2).Value = 3.684794672586015
$ws.Cells.Item(17, 3).Value = 8.953098827470686
$ws.Cells.Item(17, 4).Value = 0.5985009508893612
$ws.Cells.Item(17, 5).Value = 17
$ws.Cells.Item(17, 6).Value = 0.05343881175330965
$ws.Cells.Item(17, 7).Value = -31

$ws.Cells.Item(18, 1).Value = "Rayo Vallecano"
$ws.Cells.Item(18, 2).Value = 4.117913832199546
$ws.Cells.Item(18, 3).Value = 9.305668016194332
$ws.Cells.Item(18, 4).Value = 0.5300313122687162
$ws.Cells.Item(18, 5).Value = 18
$ws.Cells.Item(18, 6).Value = 0.0889132821075741
$ws.Cells.Item(18, 7).Value = 4

$ws.Cells.Item(19, 1).Value = "Mallorca"
$ws.Cells.Item(19, 2).Value = 5.47565543071161
$ws.Cells.Item(19, 3).Value = 7.215759849906191
$ws.Cells.Item(19, 4).Value = 0.379147465437788
$ws.Cells.Item(19, 5).Value = 12
$ws.Cells.Item(19, 6).Value = 0.1079368608270342
$ws.Cells.Item(19, 7).Value = -13

$ws.Cells.Item(20, 1).Value = "Levante UD"
$ws.Cells.Item(20, 2).Value = 5.454225352112676
$ws.Cells.Item(20, 3).Value = 6.092783505154639
$ws.Cells.Item(20, 4).Value = 0.3503397158739963
$ws.Cells.Item(20, 5).Value = 11
$ws.Cells.Item(20, 6).Value = 0.1046240647614375
$ws.Cells.Item(20, 7).Value = -7

$ws.Cells.Item(21, 1).Value = "Real Oviedo"
$ws.Cells.Item(21, 2).Value = 5.368181818181818
$ws.Cells.Item(21, 3).Value = 8.404296875
$ws.Cells.Item(21, 4).Value = 0.3644399545397146
$ws.Cells.Item(21, 5).Value = 24
$ws.Cells.Item(21, 6).Value = 0.09850034083162917
$ws.Cells.Item(21, 7).Value = -28

